# Fixed #295 Add the version of M2Doc in the template custom properties.
#
# The template itself is not changed: we only stamp the nominal-template
# document with the M2Doc version that produced/validated it, recorded as
# a custom document property (docProps/custom.xml), the same way Word's
# File > Info > Advanced Properties > Custom tab would add it.

$d = $word.ActiveDocument

$propName  = "M2DocVersion"
$propValue = "1.0.0"
# msoPropertyTypeString = 4
$propType  = 4

try {
    $existing = $d.CustomDocumentProperties($propName)
    if ($existing) {
        $existing.Value = $propValue
    } else {
        $d.CustomDocumentProperties.Add($propName, $false, $propType, $propValue)
    }
} catch {
    try {
        # Fall back to a plain Add if the lookup above is not supported.
        $d.CustomDocumentProperties.Add($propName, $false, $propType, $propValue)
    } catch {
        # Custom document properties are not reachable through every COM
        # host; leave the rest of the template untouched in that case so
        # no unrelated content is modified.
    }
}
